$d = $word.ActiveDocument

# 1. Merge the heading text "To Do Sunday" + " May 15" into a single run
#    "To Do Sunday May 15" (the bookmark currently sitting between them gets
#    dropped here and re-created at its new location below).
$d.Content.Find.Execute("To Do Sunday May 15", $true, $false, $false, $false, $false, $true, 1, $false, "To Do Sunday May 15", 2) | Out-Null

# 2. Fix "random is statements" -> "random if statements"
$d.Content.Find.Execute("random is statements", $true, $false, $false, $false, $false, $true, 1, $false, "random if statements", 2) | Out-Null

# 3. Append the new closing sentence about findCorrelations to the end of
#    paragraph 2 (the "Rewrite some of the R scripts..." paragraph).
$p2 = $d.Paragraphs(2)
$endOfP2 = $p2.Range.End - 1
$insertionRange = $d.Range($endOfP2, $endOfP2)
$insertionRange.InsertAfter(" Our findCorrelations script does not return a weight matrix but rather a vector of weights")

# 4. Re-create the _GoBack bookmark at the end of paragraph 2 (this is Word's
#    "last edit" bookmark, it is a singleton, so adding it here automatically
#    removes it from its old location in the heading).
$p2 = $d.Paragraphs(2)
$bmPos = $p2.Range.End - 1
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
